$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.805
$ws.Range("C9").Value = -11.445
$ws.Range("D11").Value = -8.316999999999998
$ws.Range("C18").Value = -12.314
$ws.Range("C20").Value = -12.581
